$d = $word.ActiveDocument

# ---------- Change 1 ----------
# "O sistema" | " exibe dados do livro selecionado e a mensagem "Deseja..."
# becomes
# "O sistema exibe a mensagem " | " "Deseja..."
$r1 = $d.Content
$found1 = $r1.Find.Execute("O sistema", $true, $false, $false, $false, $false, $true, 1, $false, `
    "O sistema exibe a mensagem ", 1)

$r2 = $d.Content
$old2 = " exibe dados do livro selecionado e a mensagem " + [char]0x201C + "Deseja realmente remover o"
$new2 = " " + [char]0x201C + "Deseja realmente remover o"
$found2 = $r2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 1)

# ---------- Change 2 ----------
$r3 = $d.Content
$found3 = $r3.Find.Execute("Tela21_Organização - excluir cupons", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Tela16_Organização - Exclusão", 1)

# ---------- Change 3 ----------
$r4 = $d.Content
$old4 = " exclui os dados do livro do meio persistente."
$new4 = " exclui os dados do cupom do meio persistente."
$found4 = $r4.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 1)

Write-Host "found1=$found1 found2=$found2 found3=$found3 found4=$found4"
